$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Variables")
$ws2 = $wb.Worksheets.Item("Categories")

# --- Sheet "Variables": fill in rows 2-8 (column B = name, column C = label, column D = valueType) ---

# Rows 2-5: short names, vertical-top alignment, black font
$ws1.Range("B2").Value = "v34_lbl55"
$ws1.Range("B2").Font.Color = 0
$ws1.Range("B2").VerticalAlignment = -4160

$ws1.Range("B3").Value = "v34_lbl54"
$ws1.Range("B3").Font.Color = 0
$ws1.Range("B3").VerticalAlignment = -4160

$ws1.Range("B4").Value = "v34_lbl57"
$ws1.Range("B4").Font.Color = 0
$ws1.Range("B4").VerticalAlignment = -4160

$ws1.Range("B5").Value = "v34_lbl56"
$ws1.Range("B5").Font.Color = 0
$ws1.Range("B5").VerticalAlignment = -4160

# Rows 6-8: longer names, wrap text, black font (explicit Calibri, no theme scheme)
$ws1.Range("B6").Value = "ant_bmi_kgm2"
$ws1.Range("B6").Font.Color = 0
$ws1.Range("B6").Font.Name = "Calibri"
$ws1.Range("B6").WrapText = $true
$ws1.Rows.Item(6).RowHeight = 29

$ws1.Range("B7").Value = "waist_cm"
$ws1.Range("B7").Font.Color = 0
$ws1.Range("B7").Font.Name = "Calibri"
$ws1.Range("B7").WrapText = $true

$ws1.Range("B8").Value = "hip_cm"
$ws1.Range("B8").Font.Color = 0
$ws1.Range("B8").Font.Name = "Calibri"
$ws1.Range("B8").WrapText = $true

# Column C (label): "?" in red, for rows 2-8
$ws1.Range("C2").Value = "?"
$ws1.Range("C2").Font.Color = 255
$ws1.Range("C3").Value = "?"
$ws1.Range("C3").Font.Color = 255
$ws1.Range("C4").Value = "?"
$ws1.Range("C4").Font.Color = 255
$ws1.Range("C5").Value = "?"
$ws1.Range("C5").Font.Color = 255
$ws1.Range("C6").Value = "?"
$ws1.Range("C6").Font.Color = 255
$ws1.Range("C7").Value = "?"
$ws1.Range("C7").Font.Color = 255
$ws1.Range("C8").Value = "?"
$ws1.Range("C8").Font.Color = 255

# Column D (valueType): "decimal" for rows 2-8, default formatting
$ws1.Range("D2").Value = "decimal"
$ws1.Range("D3").Value = "decimal"
$ws1.Range("D4").Value = "decimal"
$ws1.Range("D5").Value = "decimal"
$ws1.Range("D6").Value = "decimal"
$ws1.Range("D7").Value = "decimal"
$ws1.Range("D8").Value = "decimal"

# Restore the originally-authored selection/active cell
$ws1.Range("C12").Select()

# --- Sheet "Categories": page setup (printable) ---
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1
